# Insert two new rows at 528 (pushing the existing 528-627 down to 530-629)
# and populate the new rows with fresh data, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("528:529").Insert()

# --- New row 528 ---
$ws.Range("A528").Value = 7
$ws.Range("B528").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C528").Value = "Ñuble"
$ws.Range("D528").Value = 44637
$ws.Range("E528").Value = 16
$ws.Range("F528").Value = "Fruta"
$ws.Range("G528").Value = 100102
$ws.Range("H528").Value = "Cítricos"
$ws.Range("I528").Value = 100102003
$ws.Range("J528").Value = "Limón"
$ws.Range("K528").Value = "Sin especificar"
$ws.Range("L528").Value = "2a amarillo"
$ws.Range("M528").Value = 60
$ws.Range("N528").Value = 20000
$ws.Range("O528").Value = 20000
$ws.Range("P528").Value = 20000
$ws.Range("Q528").Value = "$/malla 16 kilos"
$ws.Range("R528").Value = "Región de O'Higgins"
$ws.Range("S528").Value = 1250
$ws.Range("T528").Value = 16

# --- New row 529 ---
$ws.Range("A529").Value = 7
$ws.Range("B529").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C529").Value = "Ñuble"
$ws.Range("D529").Value = 44637
$ws.Range("E529").Value = 16
$ws.Range("F529").Value = "Fruta"
$ws.Range("G529").Value = 100102
$ws.Range("H529").Value = "Cítricos"
$ws.Range("I529").Value = 100102003
$ws.Range("J529").Value = "Limón"
$ws.Range("K529").Value = "Sin especificar"
$ws.Range("L529").Value = "2a plateado"
$ws.Range("M529").Value = 120
$ws.Range("N529").Value = 20000
$ws.Range("O529").Value = 21000
$ws.Range("P529").Value = 20500
$ws.Range("Q529").Value = "$/malla 16 kilos"
$ws.Range("R529").Value = "Región de O'Higgins"
$ws.Range("S529").Value = 1281
$ws.Range("T529").Value = 16
